# Auto-generated edit script: update column F ("想去人数") values
# across sheets 展览, 演出, 全部类型 per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 239
$ws.Cells.Item(3, 6).Value = 74
$ws.Cells.Item(4, 6).Value = 9310
$ws.Cells.Item(5, 6).Value = 573
$ws.Cells.Item(8, 6).Value = 252
$ws.Cells.Item(9, 6).Value = 339
$ws.Cells.Item(10, 6).Value = 397
$ws.Cells.Item(12, 6).Value = 161
$ws.Cells.Item(14, 6).Value = 426
$ws.Cells.Item(15, 6).Value = 11904
$ws.Cells.Item(19, 6).Value = 152
$ws.Cells.Item(24, 6).Value = 151
$ws.Cells.Item(25, 6).Value = 2712
$ws.Cells.Item(27, 6).Value = 62
$ws.Cells.Item(29, 6).Value = 51
$ws.Cells.Item(31, 6).Value = 984
$ws.Cells.Item(32, 6).Value = 4179
$ws.Cells.Item(33, 6).Value = 3603
$ws.Cells.Item(34, 6).Value = 436
$ws.Cells.Item(35, 6).Value = 2615
$ws.Cells.Item(36, 6).Value = 3051
$ws.Cells.Item(37, 6).Value = 10
$ws.Cells.Item(38, 6).Value = 1307
$ws.Cells.Item(40, 6).Value = 770
$ws.Cells.Item(41, 6).Value = 94
$ws.Cells.Item(42, 6).Value = 409
$ws.Cells.Item(43, 6).Value = 483
$ws.Cells.Item(45, 6).Value = 131
$ws.Cells.Item(46, 6).Value = 210
$ws.Cells.Item(47, 6).Value = 107
$ws.Cells.Item(48, 6).Value = 122
$ws.Cells.Item(49, 6).Value = 127

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 14
$ws.Cells.Item(18, 6).Value = 7
$ws.Cells.Item(20, 6).Value = 6
$ws.Cells.Item(22, 6).Value = 75
$ws.Cells.Item(24, 6).Value = 34

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 239
$ws.Cells.Item(6, 6).Value = 74
$ws.Cells.Item(7, 6).Value = 9310
$ws.Cells.Item(8, 6).Value = 573
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 252
$ws.Cells.Item(12, 6).Value = 339
$ws.Cells.Item(13, 6).Value = 397
$ws.Cells.Item(14, 6).Value = 161
$ws.Cells.Item(15, 6).Value = 426
$ws.Cells.Item(16, 6).Value = 11904
$ws.Cells.Item(24, 6).Value = 151
$ws.Cells.Item(25, 6).Value = 2713
$ws.Cells.Item(27, 6).Value = 62
$ws.Cells.Item(29, 6).Value = 51
$ws.Cells.Item(30, 6).Value = 7
$ws.Cells.Item(32, 6).Value = 984
$ws.Cells.Item(33, 6).Value = 4179
$ws.Cells.Item(34, 6).Value = 3603
$ws.Cells.Item(35, 6).Value = 436
$ws.Cells.Item(36, 6).Value = 2615
$ws.Cells.Item(37, 6).Value = 3051
$ws.Cells.Item(38, 6).Value = 10
$ws.Cells.Item(39, 6).Value = 1307
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 409
$ws.Cells.Item(43, 6).Value = 483
$ws.Cells.Item(45, 6).Value = 131
$ws.Cells.Item(46, 6).Value = 210
$ws.Cells.Item(47, 6).Value = 107
$ws.Cells.Item(48, 6).Value = 122
$ws.Cells.Item(49, 6).Value = 127
